# Updates the "Price" (D) and "Volume(1h)" (E) columns on the active sheet
# to reflect the latest refreshed symbol list values. Values are written as
# literal text (matching the source data's inline-string cells), so each
# target cell is first switched to a text number format before the new
# string is assigned - otherwise Excel would auto-convert the numeric- or
# percent-looking text into a real number and silently drop things like
# trailing zeros (e.g. "297.50" -> 297.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "297.50" },
    @{ Cell = "D3";  Value = "31.65" },
    @{ Cell = "E3";  Value = "0.70%" },
    @{ Cell = "D4";  Value = "5.082" },
    @{ Cell = "D5";  Value = "0.08027" },
    @{ Cell = "E5";  Value = "8.80%" },
    @{ Cell = "D6";  Value = "2.578" },
    @{ Cell = "E6";  Value = "12.66%" },
    @{ Cell = "D7";  Value = "7.798" },
    @{ Cell = "E7";  Value = "-1.28%" },
    @{ Cell = "D9";  Value = "0.9248" },
    @{ Cell = "E9";  Value = "-0.32%" },
    @{ Cell = "D10"; Value = "0.1760" },
    @{ Cell = "E10"; Value = "2.79%" },
    @{ Cell = "D11"; Value = "0.07399" },
    @{ Cell = "E11"; Value = "-1.04%" },
    @{ Cell = "D12"; Value = "0.08996" },
    @{ Cell = "E12"; Value = "10.30%" },
    @{ Cell = "D13"; Value = "0.03057" },
    @{ Cell = "E13"; Value = "0.60%" },
    @{ Cell = "D14"; Value = "0.1001" },
    @{ Cell = "E14"; Value = "0.82%" },
    @{ Cell = "D15"; Value = "0.001504" },
    @{ Cell = "E15"; Value = "0.04%" },
    @{ Cell = "D16"; Value = "0.006038" },
    @{ Cell = "E16"; Value = "-0.71%" },
    @{ Cell = "E17"; Value = "2.65%" },
    @{ Cell = "E18"; Value = "1.23%" },
    @{ Cell = "D20"; Value = "0.1337" },
    @{ Cell = "E20"; Value = "-0.02%" },
    @{ Cell = "D21"; Value = "4.014" },
    @{ Cell = "E21"; Value = "-13.76%" },
    @{ Cell = "E22"; Value = "4.06%" },
    @{ Cell = "D23"; Value = "0.04589" },
    @{ Cell = "E23"; Value = "-1.29%" },
    @{ Cell = "D24"; Value = "0.001242" },
    @{ Cell = "E24"; Value = "1.91%" },
    @{ Cell = "E25"; Value = "-0.97%" },
    @{ Cell = "D26"; Value = "0.0001198" },
    @{ Cell = "E26"; Value = "-7.73%" },
    @{ Cell = "D27"; Value = "0.0003412" },
    @{ Cell = "E27"; Value = "82.28%" },
    @{ Cell = "D39"; Value = "0.01760" },
    @{ Cell = "E39"; Value = "2.37%" },
    @{ Cell = "D40"; Value = "0.04506" },
    @{ Cell = "E40"; Value = "-0.25%" },
    @{ Cell = "D41"; Value = "0.006847" },
    @{ Cell = "E41"; Value = "-3.96%" },
    @{ Cell = "E42"; Value = "0.23%" },
    @{ Cell = "D43"; Value = "0.002207" },
    @{ Cell = "E43"; Value = "-3.11%" },
    @{ Cell = "D44"; Value = "0.009842" },
    @{ Cell = "E44"; Value = "-6.56%" },
    @{ Cell = "D45"; Value = "0.00006459" },
    @{ Cell = "E45"; Value = "2.72%" },
    @{ Cell = "E46"; Value = "-0.03%" },
    @{ Cell = "E47"; Value = "-55.62%" },
    @{ Cell = "D48"; Value = "0.008747" },
    @{ Cell = "D49"; Value = "0.00002100" },
    @{ Cell = "E49"; Value = "-0.03%" },
    @{ Cell = "D50"; Value = "0.0002000" },
    @{ Cell = "E50"; Value = "0.04%" }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}
